$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sulphathiazole): Prediction method Use Eq. 2 -> Use Eq. 1
$ws.Range("G2").Value = "Use Eq. 1"

# Row 4 (Sulphacetamide): type acid -> base
$ws.Range("C4").Value = "base"

# Row 6 (Aspirin): Prediction method Use Eq. 1 -> Use Eq. 2
$ws.Range("G6").Value = "Use Eq. 2"

# Row 7 (Salicylic acid): Prediction method Use Eq. 1 -> Use Eq. 2
$ws.Range("G7").Value = "Use Eq. 2"

# Update the active cell selection from M13 to G13
$ws.Range("G13").Select()
